$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regs: this year's fluke minimum size (col F) is now driven off last
# year's fluke minimum size (col R) via a "minus 1" house conservation rule,
# instead of being a hard-coded number. Apply the "0.0" display format to
# both columns (R already carried half-inch values like 16.5) and wire up
# F = R - 1 as a real formula for every data row.
$ws.Range("F2:F124").NumberFormat = "0.0"
$ws.Range("R2:R124").NumberFormat = "0.0"

$ws.Range("F2").Formula = "=R2-1"
$ws.Range("F3:F66").Formula = "=R3-1"
$ws.Range("F67:F124").Formula = "=R67-1"

# Leave a selection near the top of the sheet, matching where the author
# ended up after wiring the formulas (rather than the old scrolled-down
# F109:F124 selection).
[void]$ws.Range("S15").Select()

Write-Host "done"
